$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GREETINGS")

# New row 5: GREETING_4 / Howdy. (not yet translated, excluded from word count)
$ws.Range("A5").Value = 44021
$ws.Range("A5").NumberFormat = "d-mmm"
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "GREETING_4"
$ws.Range("D5").Value = "Howdy."
$ws.Range("F5").Value = "Not translated yet, so don't include in word count."

# New row 6: GREETING_5 / HELLO! (translation identical to English, excluded from word count)
$ws.Range("A6").Value = 44021
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = "Yes"
$ws.Range("C6").Value = "GREETING_5"
$ws.Range("D6").Value = "HELLO!"
$ws.Range("E6").Value = "HELLO!"
$ws.Range("F6").Value = "Same content, so don't include in word count."

# GREETINGS becomes the active/selected sheet (was FAREWELLS before)
$ws.Activate()
$ws.Range("F7").Select()
